# Add "duration" (F) and "when to apply" (G) columns to the Universkin
# products sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (bold style, matching the other headers in row 1) ----
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "duration"
$ws.Range("G1").Value = "when to apply"

# ---- Data rows 2-9: plain (unstyled) values ----
$ws.Range("F2").Value = "60 days"
$ws.Range("G2").Value = "AM&PM"

$ws.Range("F3").Value = "60 days"
$ws.Range("G3").Value = "AM&PM"

$ws.Range("F4").Value = "30 days"
$ws.Range("G4").Value = "AM&PM"

$ws.Range("F5").Value = "60 days"
$ws.Range("G5").Value = "PM"

$ws.Range("F6").Value = "30 days"
$ws.Range("G6").Value = "AM&PM"

$ws.Range("F7").Value = "30 days"
$ws.Range("G7").Value = "AM&PM"

$ws.Range("F8").Value = "60 days"
$ws.Range("G8").Value = "AM&PM"

$ws.Range("F9").Value = "60 days"
$ws.Range("G9").Value = "AM&PM"

# ---- Data rows 10-13: column F carries the "Calibri" style used by the
# ---- rest of the Kit/Sunscreen rows (same as column D/B there) ----
$ws.Range("D10").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F10").Value = "30 days"
$ws.Range("G10").Value = "AM&PM"

$ws.Range("D11").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F11").Value = "30 days"
$ws.Range("G11").Value = "AM&PM"

$ws.Range("D12").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F12").Value = "30 days"
$ws.Range("G12").Value = "AM&PM"

$ws.Range("D10").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F13").Value = "90 days"
$ws.Range("G13").Value = "AM"

# ---- Selection, matching the post-edit workbook's active cell ----
$ws.Range("K13").Select()
